$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 37)
$ws.Range("D2").Value = [double]"1.277206020440273E-09"
$ws.Range("E2").Value = [double]"1.277206020440273E-09"

# Row 3 (Control 4)
$ws.Range("D3").Value = [double]"1.078473582667919E-36"
$ws.Range("E3").Value = [double]"1.078473582667919E-36"

# Row 4 (Control 45)
$ws.Range("D4").Value = [double]"3.77913288638544E-46"
$ws.Range("E4").Value = [double]"3.77913288638544E-46"

# Row 5 (Control 48)
$ws.Range("D5").Value = [double]"9.709866522248019E-88"
$ws.Range("E5").Value = [double]"9.709866522248019E-88"

# Row 6 (Control 20)
$ws.Range("D6").Value = [double]"3.112912560137823E-25"
$ws.Range("E6").Value = [double]"3.112912560137823E-25"

# Row 8 (MDD 12)
$ws.Range("D8").Value = [double]"0.9999999999816445"
$ws.Range("E8").Value = [double]"1.835553931073264E-11"

# Row 9 (MDD 53)
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0

# Row 10 (MDD 29)
$ws.Range("D10").Value = [double]"5.292404928633002E-18"
$ws.Range("E10").Value = 1

# Row 11 (MDD 55)
$ws.Range("D11").Value = [double]"1.685405785387158E-70"
$ws.Range("F11").Value = [double]"20.04392051696777"
